$wb = $excel.ActiveWorkbook

# Sheet "OFF" - row 3 (A3 = "R")
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 216
$wsOff.Range("C3").Value = 145
$wsOff.Range("D3").Value = 37
$wsOff.Range("E3").Value = 13
$wsOff.Range("F3").Value = 8

# Sheet "DEF" - row 3 (A3 = "R")
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 243
$wsDef.Range("C3").Value = 168
$wsDef.Range("D3").Value = 55
$wsDef.Range("E3").Value = 25
$wsDef.Range("G3").Value = 8
